$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H40").Value = 936
$ws.Range("I40").Value = 734.2857
$ws.Range("J40").Value = 1112.5
$ws.Range("K40").Value = 734.2857
$ws.Range("L40").Value = 1112.5
$ws.Range("M40").Value = -559.2857
$ws.Range("N40").Value = -1462.5

$ws.Range("H62").Value = 8710.556
$ws.Range("I62").Value = 7679.8
$ws.Range("J62").Value = 9999
$ws.Range("K62").Value = 7679.8
$ws.Range("L62").Value = 9999
$ws.Range("M62").Value = -7055.8
$ws.Range("N62").Value = -11247

$ws.Range("H65").Value = 8710.556
$ws.Range("I65").Value = 7679.8
$ws.Range("J65").Value = 9999
$ws.Range("K65").Value = 38399
$ws.Range("L65").Value = 49995
$ws.Range("M65").Value = -35279
$ws.Range("N65").Value = -56235

$ws.Range("H106").Value = 3459.8845
$ws.Range("I106").Value = 3348.6
$ws.Range("J106").Value = 3486.3809
$ws.Range("K106").Value = 3348.6
$ws.Range("L106").Value = 3486.3809
$ws.Range("M106").Value = -2717.6
$ws.Range("N106").Value = -4748.3809

$ws.Range("H111").Value = 2791.8125
$ws.Range("I111").Value = 3059.7273
$ws.Range("J111").Value = 2202.4
$ws.Range("K111").Value = 9179.1819
$ws.Range("L111").Value = 6607.200000000001
$ws.Range("M111").Value = -6112.1819
$ws.Range("N111").Value = -12741.2

$ws.Range("H112").Value = 4033341.2
$ws.Range("I112").Value = 62500500
$ws.Range("J112").Value = 1123.4482
$ws.Range("K112").Value = 187501500
$ws.Range("L112").Value = 3370.3446
$ws.Range("M112").Value = -187500392
$ws.Range("N112").Value = -5586.3446

$ws.Range("H113").Value = 58827260
$ws.Range("I113").Value = 142858820
$ws.Range("J113").Value = 5169.5
$ws.Range("K113").Value = 142858820
$ws.Range("L113").Value = 5169.5
$ws.Range("M113").Value = -142855566
$ws.Range("N113").Value = -11677.5

$ws.Range("H129").Value = 295058.22
$ws.Range("I129").Value = 265.66666
$ws.Range("J129").Value = 323586.53
$ws.Range("K129").Value = 796.9999799999999
$ws.Range("L129").Value = 970759.5900000001
$ws.Range("M129").Value = 4203.00002
$ws.Range("N129").Value = -980759.5900000001

$ws.Range("H138").Value = 4043.5806
$ws.Range("I138").Value = 3029.5454
$ws.Range("J138").Value = 4601.3
$ws.Range("K138").Value = 9088.636200000001
$ws.Range("L138").Value = 13803.9
$ws.Range("M138").Value = -3948.636200000001
$ws.Range("N138").Value = -24083.9

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H5").Value = 54.5
$ws.Range("I5").Value = 54.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 54.5
$ws.Range("L5").Value = 0
$ws.Range("N5").Value = 57.5
$ws.Range("M5").ClearContents()

$ws.Range("H32").Value = 10569.92
$ws.Range("I32").Value = 7896.2295
$ws.Range("J32").Value = 22219.572
$ws.Range("K32").Value = 7896.2295
$ws.Range("L32").Value = 22219.572
$ws.Range("M32").Value = -7609.2295
$ws.Range("N32").Value = -22793.572

$ws.Range("H61").Value = 11908915
$ws.Range("I61").Value = 13892842
$ws.Range("J61").Value = 5349.25
$ws.Range("K61").Value = 13892842
$ws.Range("L61").Value = 5349.25
$ws.Range("M61").Value = -13892630
$ws.Range("N61").Value = -5773.25

$ws.Range("H74").Value = 33335030
$ws.Range("I74").Value = 45455270
$ws.Range("J74").Value = 4362.5
$ws.Range("K74").Value = 45455270
$ws.Range("L74").Value = 4362.5
$ws.Range("M74").Value = -45454396
$ws.Range("N74").Value = -6110.5

$ws.Range("H77").Value = 33335030
$ws.Range("I77").Value = 45455270
$ws.Range("J77").Value = 4362.5
$ws.Range("K77").Value = 227276350
$ws.Range("L77").Value = 21812.5
$ws.Range("M77").Value = -227271982
$ws.Range("N77").Value = -30548.5

$ws.Range("H110").Value = 1436.5
$ws.Range("I110").Value = 1052.75
$ws.Range("J110").Value = 4506.5
$ws.Range("K110").Value = 1052.75
$ws.Range("L110").Value = 4506.5
$ws.Range("M110").Value = 992.25
$ws.Range("N110").Value = -8596.5

$ws.Range("H136").Value = 11908915
$ws.Range("I136").Value = 13892842
$ws.Range("J136").Value = 5349.25
$ws.Range("K136").Value = 41678526
$ws.Range("L136").Value = 16047.75
$ws.Range("M136").Value = -41675976
$ws.Range("N136").Value = -21147.75

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H4").Value = 54.5
$ws.Range("I4").Value = 54.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 54.5
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = 60.5
$ws.Range("M4").ClearContents()

$ws.Range("H26").Value = 3400
$ws.Range("I26").Value = 3400
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 3400
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = -3108
$ws.Range("M26").ClearContents()

$ws.Range("H36").Value = 677.3333
$ws.Range("I36").Value = 677.3333
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 677.3333
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -143.3333

$ws.Range("H82").Value = 14553
$ws.Range("I82").Value = 5994.2
$ws.Range("J82").Value = 35950
$ws.Range("K82").Value = 5994.2
$ws.Range("L82").Value = 35950
$ws.Range("M82").Value = -5611.2
$ws.Range("N82").Value = -36716

$ws.Range("H85").Value = 14553
$ws.Range("I85").Value = 5994.2
$ws.Range("J85").Value = 35950
$ws.Range("K85").Value = 5994.2
$ws.Range("L85").Value = 35950
$ws.Range("M85").Value = -4668.2
$ws.Range("N85").Value = -38602

$ws.Range("H97").Value = 12793.23
$ws.Range("I97").Value = 4758.857
$ws.Range("J97").Value = 22166.666
$ws.Range("K97").Value = 4758.857
$ws.Range("L97").Value = 22166.666
$ws.Range("M97").Value = -3767.857
$ws.Range("N97").Value = -24148.666

$ws.Range("H107").Value = 2425.0264
$ws.Range("I107").Value = 2084.3462
$ws.Range("J107").Value = 3163.1667
$ws.Range("K107").Value = 2084.3462
$ws.Range("L107").Value = 3163.1667
$ws.Range("M107").Value = -164.3462
$ws.Range("N107").Value = -7003.1667

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H16").Value = 975
$ws.Range("I16").Value = 1005.55554
$ws.Range("J16").Value = 700
$ws.Range("K16").Value = 1005.55554
$ws.Range("L16").Value = 700
$ws.Range("M16").Value = -718.55554
$ws.Range("N16").Value = -1274

$ws.Range("H31").Value = 4390.2324
$ws.Range("I31").Value = 1402.619
$ws.Range("J31").Value = 7242.0454
$ws.Range("K31").Value = 1402.619
$ws.Range("L31").Value = 7242.0454
$ws.Range("M31").Value = -1107.619
$ws.Range("N31").Value = -7832.0454

$ws.Range("H34").Value = 4390.2324
$ws.Range("I34").Value = 1402.619
$ws.Range("J34").Value = 7242.0454
$ws.Range("K34").Value = 1402.619
$ws.Range("L34").Value = 7242.0454
$ws.Range("M34").Value = -1200.619
$ws.Range("N34").Value = -7646.0454

$ws.Range("H113").Value = 975
$ws.Range("I113").Value = 1005.55554
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 1005.55554
$ws.Range("L113").Value = 700
$ws.Range("M113").Value = 1164.44446
$ws.Range("N113").Value = -5040

$ws.Range("H134").Value = 37037816
$ws.Range("I134").Value = 41667390
$ws.Range("J134").Value = 1176
$ws.Range("K134").Value = 125002170
$ws.Range("L134").Value = 3528
$ws.Range("M134").Value = -124999635
$ws.Range("N134").Value = -8598

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H24").Value = 1666
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1666
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 4998
$ws.Range("N24").Value = -5458

$ws.Range("H75").Value = 4725
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 4725
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 14175
$ws.Range("N75").Value = -16171

$ws.Range("H78").Value = 4725
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 4725
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 42525
$ws.Range("N78").Value = -52509

$ws.Range("H87").Value = 22344.5
$ws.Range("I87").Value = 9580
$ws.Range("J87").Value = 28146.545
$ws.Range("K87").Value = 28740
$ws.Range("L87").Value = 84439.63499999999
$ws.Range("M87").Value = -27492
$ws.Range("N87").Value = -86935.63499999999

$ws.Range("H90").Value = 22344.5
$ws.Range("I90").Value = 9580
$ws.Range("J90").Value = 28146.545
$ws.Range("K90").Value = 86220
$ws.Range("L90").Value = 253318.905
$ws.Range("M90").Value = -79980
$ws.Range("N90").Value = -265798.905

$ws.Range("H123").Value = 2871.4285
$ws.Range("I123").Value = 1222
$ws.Range("J123").Value = 6995
$ws.Range("K123").Value = 3666
$ws.Range("L123").Value = 20985
$ws.Range("M123").Value = -1216
$ws.Range("N123").Value = -25885

$ws.Range("H131").Value = 710.72
$ws.Range("I131").Value = 387.41666
$ws.Range("J131").Value = 754.8068
$ws.Range("K131").Value = 1162.24998
$ws.Range("L131").Value = 2264.4204
$ws.Range("M131").Value = 3877.75002
$ws.Range("N131").Value = -12344.4204

$ws.Range("H140").Value = 3238.1765
$ws.Range("I140").Value = 1616.5555
$ws.Range("J140").Value = 5062.5
$ws.Range("K140").Value = 4849.666499999999
$ws.Range("L140").Value = 15187.5
$ws.Range("M140").Value = 330.3335000000006
$ws.Range("N140").Value = -25547.5

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H122").Value = 3380.3125
$ws.Range("I122").Value = 3593.158
$ws.Range("J122").Value = 3069.2307
$ws.Range("K122").Value = 10779.474
$ws.Range("L122").Value = 9207.6921
$ws.Range("M122").Value = -8329.474
$ws.Range("N122").Value = -14107.6921

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H55").Value = 310.88235
$ws.Range("I55").Value = 315.41666
$ws.Range("J55").Value = 300
$ws.Range("K55").Value = 315.41666
$ws.Range("L55").Value = 300
$ws.Range("M55").Value = -142.41666
$ws.Range("N55").Value = -646

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 38463650
$ws.Range("I132").Value = 71429770
$ws.Range("J132").Value = 3179.6667
$ws.Range("K132").Value = 214289310
$ws.Range("L132").Value = 9539.000100000001
$ws.Range("M132").Value = -214286780
$ws.Range("N132").Value = -14599.0001
